$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 845 (the first
# "2026/12/29" entry), shifting everything below down by two rows.
$ws.Rows("845:846").Insert()

# Row 845: 2026/02/21 (Sat), hour 22
$ws.Range("A845").Value2 = "'2026/02/21"
$ws.Range("A845").ClearFormats()
$ws.Range("B845").Value2 = "土"
$ws.Range("C845").Value2 = 22
$ws.Range("D845").Value2 = 201

# Row 846: 2026/02/22 (Sun), hour 2
$ws.Range("A846").Value2 = "'2026/02/22"
$ws.Range("A846").ClearFormats()
$ws.Range("B846").Value2 = "日"
$ws.Range("C846").Value2 = 2
$ws.Range("D846").Value2 = 201
